$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current "Results" column (D),
# shifting Results to column F and leaving room for FirstName/LastName.
$ws.Range("D1:E1").EntireColumn.Insert()

# New header row additions
$ws.Range("D1").Value = "FirstName"
$ws.Range("E1").Value = "LastName"

# New data row additions
$ws.Range("D2").Value = "Test User"
$ws.Range("E2").Value = "Thomas"
$ws.Range("F2").ClearContents()

# New Browser column
$ws.Range("G1").Value = "Browser"
$ws.Range("G1").Borders.Item(7).LineStyle = 1
$ws.Range("G1").Borders.Item(10).LineStyle = 1
$ws.Range("G2").Value = "Mozilla"

# Autofit-like column widths for the new FirstName / LastName columns
$ws.Columns.Item(4).ColumnWidth = 8.43
$ws.Columns.Item(5).ColumnWidth = 8.2

# Leave selection on the last-edited cell, matching the recorded view state
[void]$ws.Range("G2").Select()
